$d = $word.ActiveDocument

# Locate the run containing the "#Report-Date#" field text and split it:
# the original run " #Report-Date#" becomes " " and a new run
# "#Report Date Here#" (bold, sz 24, szCs 24) is appended right after it.

$found = $d.Content.Find.Execute("#Report-Date#", $false, $false, $false, $false, $false, $true, 1, $false, "#Report Date Here#", 2)
